$wb = $excel.ActiveWorkbook
$app = $excel
$app.DisplayAlerts = $false

# ------------------------------------------------------------------
# 1) The "prevalenceControl" sheet needs a new internal sheetId (5 -> 6).
#    The only way to obtain a fresh sheetId in-place is to duplicate the
#    sheet (the copy receives the next free sheetId), delete the original,
#    and rename the copy back to the original name. This preserves all
#    existing formatting / the embedded table.
# ------------------------------------------------------------------
$orig = $wb.Worksheets.Item("prevalenceControl")
$orig.Copy($null, $orig)

# Re-fetch the sheet by name immediately before each mutating call:
# worksheet object handles in this runtime can become stale (they seem
# to track a positional slot) once the collection is mutated by an
# insert/delete, so never reuse a handle captured before such a change.
$original = $wb.Worksheets.Item("prevalenceControl")
$original.Delete()

$duplicate = $wb.Worksheets.Item("prevalenceControl (2)")
$duplicate.Name = "prevalenceControl"

$ws = $wb.Worksheets.Item("prevalenceControl")

# ------------------------------------------------------------------
# 2) Rename the worksheet's table from Table5 to Table1.
#    (Look the table up by its current name rather than by ordinal
#    index, since the ListObjects collection can be stale right after
#    a sheet copy/rename.)
# ------------------------------------------------------------------
$tbl = $ws.ListObjects.Item("Table5")
$tbl.Name = "Table1"

# ------------------------------------------------------------------
# 3) Update column E (Medicare Non-AS cohort) values for rows 2-24.
# ------------------------------------------------------------------
$newValues = @(
    1.3167973381,
    0.59131469940000003,
    9.6489245151999992,
    2.3395093952999999,
    2.2403606433999999,
    0.50945870419999995,
    6.4751344796000003,
    0.71426132009999999,
    0.37718808939999998,
    0.44446179139999997,
    0.0159917342,
    0.12609075820000001,
    0.072965675199999996,
    0.0099745054,
    0.028622493700000001,
    4.7482440263000001,
    0.019461127299999999,
    0.14164881830000001,
    1.3724160475,
    2.7895005320999999,
    0.37876015810000002,
    0.14126935339999999,
    0.31891312589999998
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $newValues[$i]
}
